$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" worksheet right after "2021-Q4" (and before the
#    "总计" summary sheet). It reuses the same column layout as the other
#    quarterly sheets (基金代码 / 基金名称 / 基金规模 / 股票总仓位 / 仓位占比 /
#    持有市值(亿元) / 仓位排名).
# ---------------------------------------------------------------------------
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $prevQuarter)
$q1.Name = "2022-Q1"

# Header row (bold / bordered / centered, matching the other quarter sheets)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$q1.Range("B1:H1").Font.Bold = $true
$q1.Range("B1:H1").HorizontalAlignment = -4108
$q1.Range("B1:H1").VerticalAlignment = -4160
$q1.Range("B1:H1").Borders.LineStyle = 1

# Data row
$q1.Range("A2").Value = 0
$q1.Range("A2").Font.Bold = $true
$q1.Range("A2").HorizontalAlignment = -4108
$q1.Range("A2").VerticalAlignment = -4160
$q1.Range("A2").Borders.LineStyle = 1

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "161123"
$q1.Range("C2").Value = "易方达并购重组指数（LOF）"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "4.78"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "94.71"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "3.52"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.1683"
$q1.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new top data row for 2022-Q1 and
#    push the previously existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()

$total.Range("B2:D2").ClearFormats()
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.17

$total.Range("A2").Value = 0
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# Renumber the index column (A) sequentially for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
